# rendre le site responsive et amelioration de vitesse
#
# This script reproduces (via the Excel COM object model) the changes
# described in the commit:
#  - two small text fixes on the "Feuil1" shared strings
#  - a new "Feuil2" worksheet, inserted after "Feuil1", containing an
#    expanded version of the audit table (accessibility / speed items)
#  - selection / active-sheet bookkeeping updates

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Fix a couple of typos/labels on Feuil1
# ---------------------------------------------------------------------
$ws1.Range("B2").Value = "balise meta title"
$ws1.Range("B14").Value = "definition d'une langue"

# ---------------------------------------------------------------------
# 2) Move the selection on Feuil1 from C14 to B14
# ---------------------------------------------------------------------
$ws1.Range("B14").Select()

# ---------------------------------------------------------------------
# 3) Create the new "Feuil2" worksheet right after "Feuil1".
#    Copying the header/body range from Feuil1 first means the new
#    sheet re-uses the exact same cell styles (bold header row,
#    purple fill on G1, italic/plain font on the "Reference" column,
#    etc.) instead of generating brand-new style entries.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Feuil2"

$ws1.Range("A1:F15").Copy($ws2.Range("A1"))
$ws1.Range("G1").Copy($ws2.Range("G1"))

# Drop the left-over empty cells that the block copy produced but that
# the target sheet does not have.
$ws2.Range("F6").ClearContents()
$ws2.Range("F8").ClearContents()
$ws2.Range("F9").ClearContents()
$ws2.Range("C10").ClearContents()
$ws2.Range("D10").ClearContents()
$ws2.Range("F10").ClearContents()
$ws2.Range("F11").ClearContents()
$ws2.Range("F12").ClearContents()
$ws2.Range("F13").ClearContents()
$ws2.Range("F14").ClearContents()
$ws2.Range("F15").ClearContents()

# ---------------------------------------------------------------------
# 4) Overwrite the cell values that differ from Feuil1's table so that
#    Feuil2 matches the expanded "responsive / vitesse" audit content.
# ---------------------------------------------------------------------
$ws2.Range("A2").Value = "SEO"
$ws2.Range("C2").Value = "pas de titre pour le site"

$ws2.Range("D4").Value = "Mettre des mots clés "
$ws2.Range("F4").Value = "L5/ index"

$ws2.Range("A5").Value = "SEO"
$ws2.Range("B5").Value = "Black-hat"
$ws2.Range("C5").Value = "Cacher des mots clés sur le site"
$ws2.Range("D5").Value = "Eviter les blackhat"
$ws2.Range("E5").Value = "Supprimer"
$ws2.Range("F5").Value = "L41/index-L43/index"

$ws2.Range("D6").Value = "definir des balises alt"
$ws2.Range("E6").Value = "bien décrire les images"

$ws2.Range("A7").Value = "Accessibilité"
$ws2.Range("B7").Value = "Attribut de style dans le html"
$ws2.Range("C7").Value = "Problème de maintenance"
$ws2.Range("D7").Value = "attribut dans un css"
$ws2.Range("E7").Value = "supprimer les attributs dans html"
$ws2.Range("F7").Value = "L41/index-L43/index"

$ws2.Range("A8").Value = "Accessibilité"
$ws2.Range("B8").Value = "Site pas responsive"
$ws2.Range("C8").Value = "Perte d'informations en fonction taille"
$ws2.Range("D8").Value = "Adapter le contenu du site"
$ws2.Range("E8").Value = "Rendre le site responsive"

$ws2.Range("B9").Value = "Utilisation de balise pagination"
$ws2.Range("C9").Value = "Pas de lien entre les pages"
$ws2.Range("D9").Value = "créer un lien entre les pages"
$ws2.Range("E9").Value = "ajouter pagination next"

$ws2.Range("B11").Value = "Minifier les fichiers"
$ws2.Range("C11").Value = "Alourdit nos fichiers"
$ws2.Range("D11").Value = "compresser nos css et js"

$ws2.Range("B12").Value = "Images à la place du texte"
$ws2.Range("C12").Value = "images plus lourdes que le texte "
$ws2.Range("D12").Value = "alléger le code pour un site rapide"
$ws2.Range("E12").Value = "remplacer les images par du texte"

$ws2.Range("B13").Value = "Liens defectueux"
$ws2.Range("C13").Value = "Liens cassés"
$ws2.Range("D13").Value = "Supprimer les liens morts "

$ws2.Range("B14").Value = "pas de meta robots"
$ws2.Range("C14").Value = "aucun controle d'inxation"
$ws2.Range("D14").Value = "ajouter un meta robots"

$ws2.Range("A15").Value = "SEO"
$ws2.Range("B15").Value = "definition d'une langue"
$ws2.Range("C15").Value = "site français"
$ws2.Range("D15").Value = "definir langue fr"

$ws2.Range("A16").Value = "Accessibilité"
$ws2.Range("B16").Value = "Creation des formulaires"
$ws2.Range("C16").Value = "Pas de type ou mauvaise utlisation"
$ws2.Range("D16").Value = "definir un type de formulaire"
$ws2.Range("E16").Value = "definir le bon type "

# ---------------------------------------------------------------------
# 5) Column widths for Feuil2
# ---------------------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 16.44140625
$ws2.Columns.Item(2).ColumnWidth = 31
$ws2.Columns.Item(3).ColumnWidth = 32.6640625
$ws2.Columns.Item(4).ColumnWidth = 29.6640625
$ws2.Columns.Item(5).ColumnWidth = 32.6640625
$ws2.Columns.Item(6).ColumnWidth = 31.44140625

# ---------------------------------------------------------------------
# 6) Selection / active sheet bookkeeping: Feuil2 becomes the active
#    (and selected) sheet, with B17 selected.
# ---------------------------------------------------------------------
$ws2.Range("B17").Select()
